$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (shifts existing rows 11-29 down to 12-30)
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the Natural Gas Liquids mapping
# (set C before A so new shared-string entries are created in the same
# order Excel produced them in the target workbook)
$ws.Range("C11").Value = "06_02_natural_gas_liquids"
$ws.Range("B11").Value = "06_crude_oil_and_ngl"
$ws.Range("A11").Value = "Natural Gas Liquids"

# Move the active selection to reflect the edit location
$ws.Range("D16").Select()
